$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update execution time (C2) and part (D2) for "js in browser" row
$ws.Range("C2").Value = "1h 13m"
$ws.Range("D2").Value = 2

# Update the active selection to C4
$ws.Range("C4").Select()
